# Update Name of Algo
# Applies updated KNN-imputed values to column A (and some column C) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A4"  = -20.806
    "A6"  = -22.225
    "A7"  = -19.575
    "C7"  = -12.489
    "A8"  = -22.134
    "C11" = -12.924
    "C12" = -10.782
    "C15" = -13.381
    "A16" = -21.857
    "A20" = -19.853
    "C20" = -12.366
    "A21" = -20.428
    "C21" = -12.761
    "C22" = -12.867
    "C23" = -12.223
    "A28" = -22.033
    "A29" = -21.426
    "C29" = -12.303
    "A30" = -22.112
    "A32" = -21.771
    "C34" = -11.955
    "A40" = -20.148
    "C42" = -12.492
    "C43" = -13.179
    "C44" = -13.216
    "C45" = -13.054
    "A46" = -21.951
    "C46" = -13.381
    "C50" = -13.999
    "A51" = -21.697
    "C51" = -11.2
    "A52" = -21.988
    "A57" = -21.859
    "C57" = -13.322
    "A59" = -22.278
    "A62" = -22.151
    "C65" = -12.4
    "A66" = -21.674
    "C66" = -11.405
    "C67" = -11.3
    "A73" = -20.685
    "A74" = -21.244
    "A77" = -20.342
    "C79" = -11.816
    "C84" = -13.704
    "C87" = -13.759
    "A92" = -21.626
    "C92" = -11.472
    "C97" = -12.349
    "A100" = -22.063
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
